$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = 'Components'

# Adjust window position (best-effort; xWindow offset per diff)
$win = $wb.Windows.Item(1)
$win.Left = 100

# Update component table contents
$ws.Range('A1').Value = 'Component'
$ws.Range('B1').Value = 'CType'
$ws.Range('C1').Value = 'Language'
$ws.Range('D1').Value = 'Category'
$ws.Range('E1').Value = 'Responsible'
$ws.Range('F1').Value = 'Status'
$ws.Range('G1').Value = 'Description'
$ws.Range('A2').Value = 'dataonetypes'
$ws.Range('B2').Value = 'schema'
$ws.Range('C2').Value = 'XML'
$ws.Range('D2').Value = 'Common'
$ws.Range('E2').Value = 'MJ'
$ws.Range('G2').Value = 'Schema used for defining serialization of core data types'
$ws.Range('A3').Value = 'd1_architecture'
$ws.Range('B3').Value = 'documentation'
$ws.Range('C3').Value = 'Text'
$ws.Range('D3').Value = 'Common'
$ws.Range('E3').Value = 'BS'
$ws.Range('G3').Value = 'The system architecture documentation'
$ws.Range('A4').Value = 'operations'
$ws.Range('B4').Value = 'documentation'
$ws.Range('C4').Value = 'Text'
$ws.Range('D4').Value = 'Operations'
$ws.Range('E4').Value = 'DV'
$ws.Range('G4').Value = 'Operations documentation - servers etc'
$ws.Range('A5').Value = 'd1_common_java'
$ws.Range('B5').Value = 'software'
$ws.Range('C5').Value = 'Java'
$ws.Range('D5').Value = 'Common'
$ws.Range('E5').Value = 'RW'
$ws.Range('G5').Value = 'Base DataONE library in Java'
$ws.Range('A6').Value = 'd1_common_python'
$ws.Range('B6').Value = 'software'
$ws.Range('C6').Value = 'Python'
$ws.Range('D6').Value = 'Common'
$ws.Range('E6').Value = 'RD'
$ws.Range('G6').Value = 'Base DataONE library in Python'
$ws.Range('A7').Value = 'd1_libclient_java'
$ws.Range('B7').Value = 'software'
$ws.Range('C7').Value = 'Java'
$ws.Range('D7').Value = 'Common'
$ws.Range('E7').Value = 'RN'
$ws.Range('G7').Value = 'Client library implemented in Java'
$ws.Range('A8').Value = 'd1_libclient_python'
$ws.Range('B8').Value = 'software'
$ws.Range('C8').Value = 'Python'
$ws.Range('D8').Value = 'Common'
$ws.Range('E8').Value = 'RD'
$ws.Range('G8').Value = 'Client library implemented in Python'
$ws.Range('A9').Value = 'd1_web_test_site'
$ws.Range('B9').Value = 'software'
$ws.Range('C9').Value = 'Java'
$ws.Range('D9').Value = 'Testing'
$ws.Range('E9').Value = 'RN'
$ws.Range('G9').Value = 'Member node integration testing service'
$ws.Range('A10').Value = 'd1_echo_service'
$ws.Range('B10').Value = 'software'
$ws.Range('C10').Value = 'Python'
$ws.Range('D10').Value = 'Testing'
$ws.Range('E10').Value = 'DV'
$ws.Range('G10').Value = 'A HTTP echo service used for testing'
$ws.Range('A11').Value = 'd1_integration'
$ws.Range('B11').Value = 'software'
$ws.Range('C11').Value = 'Java'
$ws.Range('D11').Value = 'Testing'
$ws.Range('E11').Value = 'RN'
$ws.Range('G11').Value = 'Integration testing for components and combinations thereof'
$ws.Range('A12').Value = 'Certificates'
$ws.Range('B12').Value = 'certificates'
$ws.Range('C12').Value = 'x509'
$ws.Range('D12').Value = 'Testing'
$ws.Range('E12').Value = 'RW'
$ws.Range('G12').Value = 'Generation and management of certificates for use by server components'
$ws.Range('A13').Value = 'd1_instance_generator'
$ws.Range('B13').Value = 'software'
$ws.Range('C13').Value = 'Python'
$ws.Range('D13').Value = 'Testing'
$ws.Range('E13').Value = 'RD'
$ws.Range('G13').Value = 'Generates example instances of objects defined in dataoneTypes.xsd'
$ws.Range('A14').Value = 'one_mercury'
$ws.Range('B14').Value = 'software'
$ws.Range('C14').Value = 'Java'
$ws.Range('D14').Value = 'CN'
$ws.Range('E14').Value = 'GP'
$ws.Range('G14').Value = 'The search interface that is implemented by the Mercury search index'
$ws.Range('A15').Value = 'cn_metacat'
$ws.Range('B15').Value = 'software'
$ws.Range('C15').Value = 'Java'
$ws.Range('D15').Value = 'CN'
$ws.Range('E15').Value = 'CJ'
$ws.Range('G15').Value = 'The `Metacat application`_. Currently employed as the replicated object  store on Coordinating Nodes.'
$ws.Range('A16').Value = 'd1_cn_index_processor'
$ws.Range('B16').Value = 'software'
$ws.Range('C16').Value = 'Java'
$ws.Range('D16').Value = 'CN'
$ws.Range('E16').Value = 'SR/DV'
$ws.Range('G16').Value = 'Populates the SOLR index by extracting informaton from system metadata, science metadata and resource maps.'
$ws.Range('A17').Value = 'd1_cn_index_generator'
$ws.Range('B17').Value = 'software'
$ws.Range('C17').Value = 'Java'
$ws.Range('D17').Value = 'CN'
$ws.Range('E17').Value = 'SR/DV'
$ws.Range('G17').Value = 'Generates indexing tasks when new objects appear or system metadata changes'
$ws.Range('A18').Value = 'd1_cn_index_common'
$ws.Range('B18').Value = 'software'
$ws.Range('C18').Value = 'Java'
$ws.Range('D18').Value = 'CN'
$ws.Range('E18').Value = 'SR/DV'
$ws.Range('G18').Value = 'Code shared between the indexing components'
$ws.Range('A19').Value = 'indexerapi'
$ws.Range('B19').Value = 'software'
$ws.Range('C19').Value = 'Java'
$ws.Range('D19').Value = 'CN'
$ws.Range('E19').Value = 'SR/DV'
$ws.Range('G19').Value = 'A library used by the index_processor for extracting content from various types of XML structures such as system metadata, science metadata and resource maps. '
$ws.Range('A20').Value = 'd1_portal_servlet'
$ws.Range('B20').Value = 'software'
$ws.Range('C20').Value = 'Java'
$ws.Range('D20').Value = 'CN'
$ws.Range('E20').Value = 'BL'
$ws.Range('G20').Value = 'Provides a UI for interacting with the CILogon service, an authentication proxy service'
$ws.Range('A21').Value = 'd1_portal'
$ws.Range('B21').Value = 'software'
$ws.Range('C21').Value = 'Java'
$ws.Range('D21').Value = 'CN'
$ws.Range('E21').Value = 'BL'
$ws.Range('G21').Value = 'Implements the certificate manager used by the portal servlets'
$ws.Range('A22').Value = 'd1_identity_manager'
$ws.Range('B22').Value = 'software'
$ws.Range('C22').Value = 'Java'
$ws.Range('D22').Value = 'CN'
$ws.Range('E22').Value = 'BL'
$ws.Range('G22').Value = 'Provides mechanisms for managing subjects in dataone'
$ws.Range('A23').Value = 'd1_process_daemon'
$ws.Range('B23').Value = 'software'
$ws.Range('C23').Value = 'Java'
$ws.Range('D23').Value = 'CN'
$ws.Range('E23').Value = 'RW'
$ws.Range('G23').Value = 'Monitors content on member nodes, creating tasks for synchronization and replication'
$ws.Range('A24').Value = 'd1_synchronization'
$ws.Range('B24').Value = 'software'
$ws.Range('C24').Value = 'Java'
$ws.Range('D24').Value = 'CN'
$ws.Range('E24').Value = 'RW'
$ws.Range('G24').Value = 'Manages the synchronization of content between Member Nodes and the Coordinating Nodes.'
$ws.Range('A25').Value = 'd1_replication'
$ws.Range('B25').Value = 'software'
$ws.Range('C25').Value = 'Java'
$ws.Range('D25').Value = 'CN'
$ws.Range('E25').Value = 'CJ'
$ws.Range('G25').Value = 'Manages replication of content between Member Nodes'
$ws.Range('A26').Value = 'd1_cn_noderegistry'
$ws.Range('B26').Value = 'software'
$ws.Range('C26').Value = 'Java'
$ws.Range('D26').Value = 'CN'
$ws.Range('E26').Value = 'RW'
$ws.Range('G26').Value = 'A register of coordinating and member nodes participating in a DataONE environment'
$ws.Range('A27').Value = 'd1_cn_common'
$ws.Range('B27').Value = 'software'
$ws.Range('C27').Value = 'Java'
$ws.Range('D27').Value = 'CN'
$ws.Range('E27').Value = 'RW'
$ws.Range('G27').Value = 'A library of code shared between coordintating node components'
$ws.Range('A28').Value = 'd1_cn_rest'
$ws.Range('B28').Value = 'software'
$ws.Range('C28').Value = 'Java'
$ws.Range('D28').Value = 'CN'
$ws.Range('E28').Value = 'RW'
$ws.Range('G28').Value = 'The coordinating node HTTP REST service interface'
$ws.Range('A29').Value = 'd1_cn_rest_proxy'
$ws.Range('B29').Value = 'software'
$ws.Range('C29').Value = 'Java'
$ws.Range('D29').Value = 'CN'
$ws.Range('E29').Value = 'RW'
$ws.Range('G29').Value = 'Proxies requests coming in to a CN to underlying service implementations such as the object store (i.e. Metacat)'
$ws.Range('A30').Value = 'd1_cn_service'
$ws.Range('B30').Value = 'software'
$ws.Range('C30').Value = 'Java'
$ws.Range('D30').Value = 'CN'
$ws.Range('E30').Value = 'RW'
$ws.Range('G30').Value = 'Coordinating node service, implementing the service APIs, data storage, and CN replication.'
$ws.Range('A31').Value = 'd1_simple_search'
$ws.Range('B31').Value = 'software'
$ws.Range('C31').Value = 'Java'
$ws.Range('D31').Value = 'CN'
$ws.Range('E31').Value = 'DV'
$ws.Range('G31').Value = 'A simple search interface using Javascript and the SOLR interface.'
$ws.Range('A32').Value = 'Metacat'
$ws.Range('B32').Value = 'software'
$ws.Range('C32').Value = 'Java'
$ws.Range('D32').Value = 'MN'
$ws.Range('E32').Value = 'CJ'
$ws.Range('G32').Value = 'The `Metacat application`_. Implements the DataONE MN service interfaces.'
$ws.Range('A33').Value = 'Dryad'
$ws.Range('B33').Value = 'software'
$ws.Range('C33').Value = 'Java'
$ws.Range('D33').Value = 'MN'
$ws.Range('E33').Value = 'RS'
$ws.Range('G33').Value = 'A member node implementation and instance for the Dryad repository'
$ws.Range('A34').Value = 'GMN'
$ws.Range('B34').Value = 'software'
$ws.Range('C34').Value = 'Python'
$ws.Range('D34').Value = 'MN'
$ws.Range('E34').Value = 'RD'
$ws.Range('G34').Value = 'A generic, standalone Member Node implementation written in Python using the Django_ framework.'
$ws.Range('A35').Value = 'Mercury_MN'
$ws.Range('B35').Value = 'software'
$ws.Range('C35').Value = 'Java'
$ws.Range('D35').Value = 'MN'
$ws.Range('E35').Value = 'JG'
$ws.Range('G35').Value = 'Mercury implementation of the Member Node services'
$ws.Range('A36').Value = 'd1_client_cli'
$ws.Range('B36').Value = 'software'
$ws.Range('C36').Value = 'Python'
$ws.Range('D36').Value = 'ITK'
$ws.Range('E36').Value = 'RD'
$ws.Range('G36').Value = 'A command line client for interacting with the DataONE infrastructure. Currently implemented using d1_libclient_python.'
$ws.Range('A37').Value = 'd1_client_fuse'
$ws.Range('B37').Value = 'software'
$ws.Range('C37').Value = 'Python'
$ws.Range('D37').Value = 'ITK'
$ws.Range('E37').Value = 'DV'
$ws.Range('G37').Value = 'A FUSE_ driver for mounting the DataONE infrastructure as a file system.'
$ws.Range('A38').Value = 'd1_client_dokan'
$ws.Range('B38').Value = 'software'
$ws.Range('C38').Value = 'Python'
$ws.Range('D38').Value = 'ITK'
$ws.Range('E38').Value = 'DV'
$ws.Range('G38').Value = 'An extention of the FUSE driver that is based on Dokan_ for use on Microsoft Windows systems.'
$ws.Range('A39').Value = 'd1_client_r'
$ws.Range('B39').Value = 'software'
$ws.Range('C39').Value = 'Java'
$ws.Range('D39').Value = 'ITK'
$ws.Range('E39').Value = 'MJ'
$ws.Range('G39').Value = 'A plugin for R that enables access to DataONE content from the R_ application. Implemented using d1_libclient_java.'
$ws.Range('A40').Value = 'hzpeek'
$ws.Range('B40').Value = 'software'
$ws.Range('C40').Value = 'Java'
$ws.Range('D40').Value = 'Testing'
$ws.Range('E40').Value = 'DV'
$ws.Range('G40').Value = 'A tool for examining the Hazelcast queues on the CNs'
$ws.Range('A41').Value = 'debian_packaing'
$ws.Range('B41').Value = 'installer'
$ws.Range('C41').Value = 'various'
$ws.Range('D41').Value = 'CN'
$ws.Range('E41').Value = 'RW'
$ws.Range('G41').Value = 'Debian packages for the CN components'

# Restore the active selection
$ws.Range('D32').Select()
